$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$r = $ws.Range("M9")
$r.Formula = "=L9/L8-1"
$r.NumberFormat = "0.0%"
